$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style used by the other
# header cells (e.g. H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I0 (column I) and IF (column J) values for rows 2-69.
$iValues = @(8,7,7,7,7,9,7,7,6,6,9,9,5,5,8,10,9,8,6,8,8,8,7,7,8,7,9,7,6,7,7,5,7,7,7,7,8,7,7,7,7,8,9,8,8,6,8,5,8,7,8,9,8,5,8,8,7,6,8,8,8,6,4,8,7,6,6,6)
$jValues = @(8,7,7,7,7,9,7,7,6,6,9,9,5,5,9,10,9,8,7,8,9,9,7,7,9,8,9,8,7,7,7,6,7,7,7,7,9,7,7,7,7,8,9,8,8,7,8,5,9,7,8,9,8,6,8,8,7,6,8,9,8,7,5,8,8,6,7,7)

for ($k = 0; $k -lt $iValues.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}
